$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that only contained "5840560 - Marco Antonio Carvalho Pereira" in
# columns B/C (with no label in column A) is removed; all subsequent rows
# shift up by one.
$ws.Rows.Item(13).Delete()

# After the shift, update the B/C values on the rows whose content changed.
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

$ws.Range("B15").Value = "01/01/2021"
$ws.Range("C15").Value = "01/01/2021"

$ws.Range("B18").Value = "5840560 - Marco Antonio Carvalho Pereira"
$ws.Range("C18").Value = "5840560 - Marco Antonio Carvalho Pereira"

$ws.Range("B19").Value = "Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras."
$ws.Range("C19").Value = "Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras."

$ws.Range("B20").Value = "Avaliação individual (Peso entre 20-40%) e do projeto realizado em equipe (peso entre 60-80%)"
$ws.Range("C20").Value = "Avaliação individual (Peso entre 20-40%) e do projeto realizado em equipe (peso entre 60-80%)"

$ws.Range("B21").Value = "NF = (MF + PR)/2, onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota darecuperação."
$ws.Range("C21").Value = "NF = (MF + PR)/2, onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota darecuperação."
